# "completed case 1 bullington method added NextPoint Excel function"
#
# Renames Sheet1 -> ProfilePoints, adds two new shared strings ("Next",
# "Intermediate Points"), replaces the old B6 "=GetDistance" formula with
# its literal cached value, and builds the new "Next point" (row 8) +
# "Intermediate Points" (rows 10-22) tables that drive the Bullington
# diffraction-loss profile, using the custom _xll.NextPoint /
# _xll.GetProfilePoints / _xll.GetDistance add-in functions. The old
# A9:B20 MakeList() demo array is removed to make room for the new
# A11:C22 profile-points array.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Sheet rename --------------------------------------------------------
$ws.Name = "ProfilePoints"

# --- Window bookkeeping (workbookView xWindow/yWindow) -------------------
$wb.Windows.Item(1).Left = 28680
$wb.Windows.Item(1).Top = 30

# --- Remove the old demo array (A9:B20 = _xll.MakeList()) ----------------
# Has to be un-arrayed over its own original extent before it can be
# touched piecemeal, otherwise Excel refuses ("You cannot change part of
# an array").
$ws.Range("A9:B20").Value = ""
$ws.Range("A9:B20").ClearContents()

# --- Row 6: Tx/Rx distance becomes a plain cached value -------------------
$ws.Range("B6").Value = 1001

# --- Row 8: "Next" point (array formula over B8:C8) -----------------------
$ws.Range("A8").Value = "Next"
$ws.Range("B8:C8").FormulaArray = "=_xll.NextPoint(B2,C2,B5,B6)"
$ws.Range("C8").Value = 151.21025912980389

# --- Row 10: section header -----------------------------------------------
$ws.Range("A10").Value = "Intermediate Points"

# --- Row 11: profile points (array formula over A11:C22) ------------------
$ws.Range("A11:C22").FormulaArray = "=_xll.GetProfilePoints(B2,C2,B8,C8)"

# --- Rows 11-22: lat / lon / elevation profile values ----------------------
$profile = @(
    @(11, -33.839535000000012, 151.20694600000002, 95),
    @(12, -33.84039321526285,  151.20727695218665, 85),
    @(13, -33.841251429517683, 151.20760791098985, 62),
    @(14, -33.842109642764484, 151.20793887640997, 45),
    @(15, -33.842967855003167, 151.20826984844734, 28),
    @(16, -33.843826066233703, 151.20860082710232, 10),
    @(17, -33.844684276456036, 151.20893181237523, 2),
    @(18, -33.845542485670137, 151.20926280426636, 2),
    @(19, -33.84640069387595,  151.20959380277611, 4),
    @(20, -33.847258901073417, 151.20992480790483, 12),
    @(21, -33.848117107262482, 151.21025581965284, 7),
    @(22, -33.848125689319261, 151.21025912980389, 7)
)

foreach ($point in $profile) {
    $r = $point[0]
    if ($r -ne 11) {
        $ws.Range("A$r").Value = $point[1]
        $ws.Range("B$r").Value = $point[2]
    }
    $ws.Range("C$r").Value = $point[3]
}

# --- Rows 12-22: per-segment distance formulas in column D ----------------
for ($r = 12; $r -le 22; $r++) {
    $prev = $r - 1
    $ws.Range("D$r").Formula = "=_xll.GetDistance(A$prev,B$prev,A$r,B$r)"
}

# --- Selection shown when the workbook is reopened -------------------------
$ws.Range("J10").Select() | Out-Null
